$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "26.528.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'" + "1.737.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").Value = "'" + "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'" + "246.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'" + "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'" + "0.4898"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").Value = "'" + "0.2668"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'" + "0.06329"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").Value = "'" + "1.730.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'" + "0.07042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "'" + "15.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'" + "4.609"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "'" + "0.6117"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'" + "77.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "'" + "0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'" + "0.000007395"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.34%  "
$ws.Range("D18").Value = "'" + "26.531.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'" + "0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'" + "1.953.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'" + "8.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("D24").Value = "'" + "5.241"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "'" + "140.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("D26").Value = "'" + "15.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'" + "1.409"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'" + "1.767"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'" + "107.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").Value = "'" + "4.032"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "'" + "0.08052"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").Value = "'" + "3.715"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "'" + "0.04579"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'" + "0.9989"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'" + "2.608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "'" + "1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("D37").Value = "'" + "0.6363"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'" + "0.8952"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.02%  "
$ws.Range("D39").Value = "'" + "2.020"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'" + "0.01505"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'" + "102.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.27%  "
$ws.Range("D44").Value = "'" + "5.396"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("D45").Value = "'" + "0.3898"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'" + "6.891"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'" + "0.1186"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'" + "30.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'" + "7.798"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").Value = "'" + "1.267"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "

Write-Output "Applied cryptos update"
